# Apply "ubah perhitungan manual agar signifikan" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the raw input "matrik keputusan" values (columns B-F, rows 4-13) ---
$ws.Range("C6").Value = 20
$ws.Range("C7").Value = 50
$ws.Range("C9").Value = 50

# --- Update the weight (bobot) inputs on row 14 ---
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 5

# Recalculate every formula-driven cell (G14/B18:G18/B22:G32/B35:B44/etc.)
$excel.Calculate()

# --- Refresh the manually-sorted "Preferensi V" helper table (D35:E44) ---
# It is a static (non-formula) copy of A35:B44 sorted descending by score,
# so it has to be rebuilt by hand after the inputs above change the scores.
$ws.Range("A35:B44").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = 0

[void]$ws.Range("D35:E44").Sort($ws.Range("E35"), 2)

# --- Move the active selection to match the saved view state ---
[void]$ws.Range("G16").Select()
